$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ITEM"
$ws.Range("B2").Value = "DESCRIÇÃO"
$ws.Range("C2").Value = "MARCA"
$ws.Range("D2").Value = "Não encontrado"
$ws.Range("E2").Value = "Não encontrado"
$ws.Range("F2").Value = "Pendente"
